# Weekly update: a new price observation is inserted as row 490 on the
# single data sheet, pushing the existing rows 490-513 down to 491-514
# (dimension grows from A1:R513 to A1:R514).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 490, shifting rows 490:513 down to 491:514.
$ws.Rows.Item(490).Insert()

# Populate the newly inserted row 490 with the new observation.
$ws.Cells.Item(490, 1).Value = 9
$ws.Cells.Item(490, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(490, 3).Value = "Metropolitana"
$ws.Cells.Item(490, 4).Value = 44939
$ws.Cells.Item(490, 5).Value = 13
$ws.Cells.Item(490, 6).Value = 100112039
$ws.Cells.Item(490, 7).Value = "Ciboulette"
$ws.Cells.Item(490, 8).Value = "Sin especificar"
$ws.Cells.Item(490, 9).Value = "Primera"
$ws.Cells.Item(490, 10).Value = 430
$ws.Cells.Item(490, 11).Value = 1000
$ws.Cells.Item(490, 12).Value = 1000
$ws.Cells.Item(490, 13).Value = 1000
$ws.Cells.Item(490, 14).Value = "`$/docena de atados"
$ws.Cells.Item(490, 15).Value = "Región Metropolitana"
$ws.Cells.Item(490, 16).Value = 333
$ws.Cells.Item(490, 17).Value = 3
$ws.Cells.Item(490, 18).Value = "Hortaliza"
